$d = $word.ActiveDocument

$pairs = @(
    @{old="17×11=187"; new="58×70=4060"},
    @{old="34×13=442"; new="82×61=5002"},
    @{old="77×76=5852"; new="30×77=2310"},
    @{old="93×21=1953"; new="93×69=6417"},
    @{old="89×54=4806"; new="87×39=3393"},
    @{old="69×26=1794"; new="38×88=3344"},
    @{old="29×25=725"; new="50×98=4900"},
    @{old="28×44=1232"; new="67×38=2546"},
    @{old="98×97=9506"; new="44×42=1848"},
    @{old="92×65=5980"; new="37×25=925"},
    @{old="45×40=1800"; new="27×97=2619"},
    @{old="93×89=8277"; new="55×58=3190"},
    @{old="67×46=3082"; new="38×26=988"},
    @{old="45×89=4005"; new="76×37=2812"},
    @{old="65×41=2665"; new="20×27=540"},
    @{old="51×84=4284"; new="88×38=3344"},
    @{old="73×65=4745"; new="18×59=1062"},
    @{old="81×99=8019"; new="34×95=3230"},
    @{old="60×94=5640"; new="43×57=2451"},
    @{old="69×51=3519"; new="82×47=3854"},
    @{old="34×70=2380"; new="76×85=6460"},
    @{old="45×61=2745"; new="71×72=5112"},
    @{old="74×13=962"; new="45×98=4410"},
    @{old="90×65=5850"; new="99×86=8514"},
    @{old="32×59=1888"; new="76×94=7144"}
)

foreach ($p in $pairs) {
    $range = $d.Content
    $range.Find.Execute($p.old, $true, $true, $false, $false, $false, $true, 1, $false, $p.new, 2)
}
